$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").HorizontalAlignment = -4108
$ws.Range("I1:J1").VerticalAlignment = -4160
$ws.Range("I1:J1").Borders.LineStyle = 1
$ws.Range("I1:J1").Borders.Weight = 2

$data = @(
    @(5, 9),
    @(1, 6),
    @(1, 6),
    @(1, 7),
    @(1, 5),
    @(1, 9),
    @(1, 4),
    @(1, 3),
    @(1, 6),
    @(1, 7),
    @(1, 7),
    @(1, 3),
    @(1, 6),
    @(1, 5),
    @(1, 9),
    @(1, 7),
    @(1, 6),
    @(1, 7),
    @(1, 4),
    @(1, 5),
    @(1, 3),
    @(1, 3),
    @(1, 3),
    @(1, 3)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
